# The commit "regenerate orders with updates distance/sizes" renames the
# trial-order's distance codes (D51->D55, D64->D69, D80->D86) and the "large"
# size code (S30->S31) everywhere they occur: in the Condition, Filename_Left,
# Filename_Right and Distance/Size columns (plain values, and as substrings of
# names like "Face08_D51_S30" / "Face08_D51_S30_l.png" / "Fixation_D51_l.png").
# None of the replacement tokens collide with each other or with pre-existing
# text, so a simple global substring Find/Replace over the used range
# reproduces the diff exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = $ws.UsedRange

# xlPart (-4163): match the token anywhere inside a cell's text, not just
# whole-cell matches (needed for e.g. "Face08_D51_S30_l.png").
$xlPart = -4163

$cells.Replace("D51", "D55", $xlPart)
$cells.Replace("D64", "D69", $xlPart)
$cells.Replace("D80", "D86", $xlPart)
$cells.Replace("S30", "S31", $xlPart)
